$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same bold/centered/bordered style as the other headers (copy format from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$times = @(
    "2021-10-05 13:41:50.395605",
    "2021-10-05 13:41:50.395617",
    "2021-10-05 13:41:50.395621",
    "2021-10-05 13:41:50.395624",
    "2021-10-05 13:41:50.395628",
    "2021-10-05 13:41:50.395631",
    "2021-10-05 13:41:50.395634",
    "2021-10-05 13:41:50.395637",
    "2021-10-05 13:41:50.395640",
    "2021-10-05 13:41:50.395643",
    "2021-10-05 13:41:50.395646",
    "2021-10-05 13:41:50.395649",
    "2021-10-05 13:41:50.395652",
    "2021-10-05 13:41:50.395655",
    "2021-10-05 13:41:50.395658",
    "2021-10-05 13:41:50.395661",
    "2021-10-05 13:41:50.395664",
    "2021-10-05 13:41:50.395667",
    "2021-10-05 13:41:50.395670",
    "2021-10-05 13:41:50.395673",
    "2021-10-05 13:41:50.395676",
    "2021-10-05 13:41:50.395679",
    "2021-10-05 13:41:50.395682",
    "2021-10-05 13:41:50.395684",
    "2021-10-05 13:41:50.395688",
    "2021-10-05 13:41:50.395691"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
